# Texas COVID data book - 2020-10-02 daily update.
# The previous commit accidentally duplicated the prior day's numbers; this
# edit replaces them with the corrected counts and turns the "%" column (and
# the "Total"/"Grand Total" row) into live formulas instead of hard-coded
# numbers, matching how the other tabs in the workbook already worked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Cases by Age Group
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value2  = 187
$ws.Range("B3").Value2  = 883
$ws.Range("B4").Value2  = 2320
$ws.Range("B5").Value2  = 10476
$ws.Range("B6").Value2  = 11559
$ws.Range("B7").Value2  = 10446
$ws.Range("B8").Value2  = 9146
$ws.Range("B9").Value2  = 3388
$ws.Range("B10").Value2 = 2304
$ws.Range("B11").Value2 = 1379
$ws.Range("B12").Value2 = 891
$ws.Range("B13").Value2 = 1420
$ws.Range("B15").Formula = '=SUM(B2:B14)'
$ws.Range("C2").Formula  = '=B2/$B$15'
$ws.Range("C3:C15").Formula = '=B3/$B$15'

# ---------------------------------------------------------------------
# 2) Cases by Gender
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value2 = 18533
$ws.Range("B3").Value2 = 34992
$ws.Range("B4").Value2 = 893
$ws.Range("B5").Formula = '=SUM(B2:B4)'
$ws.Range("C2").Formula = '=B2/$B$5'
$ws.Range("C3:C5").Formula = '=B3/$B$5'

# ---------------------------------------------------------------------
# 3) Cases by RaceEthnicity
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value2 = 812
$ws.Range("B3").Value2 = 9027
$ws.Range("B4").Value2 = 21472
$ws.Range("B6").Value2 = 16358
$ws.Range("B7").Value2 = 6507
$ws.Range("B8").Formula = '=SUM(B2:B7)'
$ws.Range("C2").Formula = '=B2/$B$8'
$ws.Range("C3:C8").Formula = '=B3/$B$8'

# ---------------------------------------------------------------------
# 4) Fatalities by Age Group
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B2").Value2  = 3
$ws.Range("B3").Value2  = 6
$ws.Range("B4").Value2  = 18
$ws.Range("B5").Value2  = 102
$ws.Range("B6").Value2  = 298
$ws.Range("B7").Value2  = 893
$ws.Range("B8").Value2  = 1952
$ws.Range("B9").Value2  = 1523
$ws.Range("B10").Value2 = 1851
$ws.Range("B11").Value2 = 2082
$ws.Range("B12").Value2 = 1950
$ws.Range("B13").Value2 = 5217
$ws.Range("B15").Formula = '=SUM(B2:B14)'
$ws.Range("C2").Formula  = '=B2/$B$15'
$ws.Range("C3:C15").Formula = '=B3/$B$15'

# ---------------------------------------------------------------------
# 5) Fatalities by Gender
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value2 = 6642
$ws.Range("B3").Value2 = 9253
$ws.Range("B5").Formula = '=SUM(B2:B4)'
$ws.Range("C2").Formula = '=B2/$B$5'
$ws.Range("C3:C5").Formula = '=B3/$B$5'

# ---------------------------------------------------------------------
# 6) Fatalities by Race-Ethnicity (note: Grand Total stays a plain number
#    here, not a formula - matches the source workbook as edited)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value2 = 295
$ws.Range("B3").Value2 = 1774
$ws.Range("B4").Value2 = 8916
$ws.Range("B5").Value2 = 87
$ws.Range("B6").Value2 = 4813
$ws.Range("B7").Value2 = 10
$ws.Range("B8").Value2 = 15895
$ws.Range("C2").Formula = '=B2/$B$8'
$ws.Range("C3:C8").Formula = '=B3/$B$8'

# ---------------------------------------------------------------------
# Selected tab moves from "Fatalities by Race-Ethnicity" back to
# "Cases by Age Group".
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Cases by Age Group").Activate()
